# Update the weekly progress log for 黃煒翔 (row 11) with the latest status,
# and reuse the new task description ("蒐集市場實際數據") on row 13 as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("工作表1")

$ws.Range("D11").Value = "草擬機率運算架構"
$ws.Range("E11").Value = "完成"
$ws.Range("F11").Value = "蒐集市場實際數據"

$ws.Range("D13").Value = "蒐集市場實際數據"

# Update the active selection to match the new edit position.
$ws.Range("F15").Select()
